# Regenerate merged AHB files
# 1) Rename the header labels from *_old / *_new to *_FV2310 / *_FV2404
# 2) Freeze the header row (split/pane)
# 3) Turn the data range into a proper Excel Table (ListObject)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A..J -> *_FV2310
$fv2310 = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

# Columns L..U -> *_FV2404
$fv2404 = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $fv2310.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2310[$i]
}

for ($i = 0; $i -lt $fv2404.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2404[$i]
}

# Freeze the header row: split after row 1, keep viewport at A2
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn A1:U64 into an Excel Table named Table1
$range = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"
